# Updates the benchmark-stats table in the Renaissance/ShenandoahGC docx.
# Rows 1-3 collapse to "0M", rows 4-12 get refreshed numeric values, and
# rows 44-46 (previously large tab-separated summary rows) collapse down
# to the single leading value that used to live in rows 1-3.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "2022"
$t.Cell(5, 1).Range.Text  = "0.00003"
$t.Cell(6, 1).Range.Text  = "0.00083"
$t.Cell(7, 1).Range.Text  = "0.00016"
$t.Cell(8, 1).Range.Text  = "0.00005"
$t.Cell(9, 1).Range.Text  = "0.00026"
$t.Cell(10, 1).Range.Text = "0.00029"
$t.Cell(11, 1).Range.Text = "0.00042"
$t.Cell(12, 1).Range.Text = "0.37708"

$t.Cell(44, 1).Range.Text = "99.82"
$t.Cell(45, 1).Range.Text = "0.38"
$t.Cell(46, 1).Range.Text = "213"
